$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price for row 342 (B342: 5800 -> 6200)
$ws.Range("B342").Value = 6200

# Append new rows 353-360
$ws.Range("A353").Value = "Букет 2286"
$ws.Range("B353").Value = 4600
$ws.Range("D353").Value = "https://gift2gift.ru/img/work/nomencl/2286-s.jpeg"
$ws.Range("E353").Value = "https://gift2gift.ru/catalog/vitrina-bukety-v-nalichii-pryamo-seichas/_______________________________________________________________________________________________________________________________________________________________________________________________________________________________________________________________.html"

$ws.Range("A354").Value = "Букет 2287"
$ws.Range("B354").Value = 4250
$ws.Range("D354").Value = "https://gift2gift.ru/img/work/nomencl/2287-s.jpeg"
$ws.Range("E354").Value = "https://gift2gift.ru/catalog/vitrina-bukety-v-nalichii-pryamo-seichas/_______________________________________________________________________________________________________________________________________________________________________________________________________________________________________________________________.html"

$ws.Range("A355").Value = "Сборный букет №111"
$ws.Range("B355").Value = 4250
$ws.Range("D355").Value = "https://gift2gift.ru/img/work/nomencl/2300-s.jpeg"
$ws.Range("E355").Value = "https://gift2gift.ru/catalog/bukety-ot-3000-do-4000/sbornyi-buket-111.html"

$ws.Range("A356").Value = "Сборный букет №112"
$ws.Range("B356").Value = 4600
$ws.Range("D356").Value = "https://gift2gift.ru/img/work/nomencl/2301-s.jpeg"
$ws.Range("E356").Value = "https://gift2gift.ru/catalog/bukety-ot-3000-do-4000/sbornyi-buket-112.html"

$ws.Range("A357").Value = "Композиция 111"
$ws.Range("B357").Value = 2730
$ws.Range("D357").Value = "https://gift2gift.ru/img/work/nomencl/2302-s.jpeg"
$ws.Range("E357").Value = "https://gift2gift.ru/catalog/cvety_v_korobke/kompozitsiya-111.html"

$ws.Range("A358").Value = "Композиция 112"
$ws.Range("B358").Value = 2545
$ws.Range("D358").Value = "https://gift2gift.ru/img/work/nomencl/2303-s.jpeg"
$ws.Range("E358").Value = "https://gift2gift.ru/catalog/cvety_v_korobke/kompozitsiya-112.html"

$ws.Range("A359").Value = "Композиция 113"
$ws.Range("B359").Value = 2370
$ws.Range("D359").Value = "https://gift2gift.ru/img/work/nomencl/2304-s.jpeg"
$ws.Range("E359").Value = "https://gift2gift.ru/catalog/cvety_v_korobke/kompozitsiya-113.html"

$ws.Range("A360").Value = "Композиция 114"
$ws.Range("B360").Value = 2170
$ws.Range("D360").Value = "https://gift2gift.ru/img/work/nomencl/2305-s.jpeg"
$ws.Range("E360").Value = "https://gift2gift.ru/catalog/cvety_v_korobke/_______________________________________________________________________________________________________________________________________________________________________________________________________________________________________________________________.html"
